# Fix mistake in excel: rows 49-66 in column B should say "OpEx" instead
# of "CapEx" (a copy/paste leftover), and update the sheet's scroll/
# selection state to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the CapEx -> OpEx mistake for rows 49 through 66 in column B.
for ($row = 49; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value() -eq "CapEx") {
        $cell.Value = "OpEx"
    }
}

# Update the view to match the state after the fix: scrolled down with
# B66 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("B66").Select()
